# Algorithm study plan sheet update:
# - Swap the "탐색" / "DFS/BFS (한 문제는 꼭 나옴)" labels between A9 and A10
# - B9 gains a " (DFS)" suffix (11724 -> "11724 (DFS)")
# - C9 gets a new date value 240522
# - Active selection on 계획표 moves to C9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("계획표")

$ws.Range("A9").Value = "DFS/BFS (한 문제는 꼭 나옴)"
$ws.Range("A10").Value = "탐색"

$ws.Range("B9").Value = "11724 (DFS)"
$ws.Range("C9").Value = 240522
$ws.Range("C9").HorizontalAlignment = -4108
$ws.Range("C9").VerticalAlignment = -4108

$ws.Range("C9").Select()
